$d = $word.ActiveDocument

# Locate the paragraph that ends with "LOB1049: Estatística Multivariada (Requisito fraco)".
# The three paragraphs that immediately follow it (an empty paragraph, the
# "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph, and the
# copyright/footer paragraph) must be removed, while the paragraph after
# those (another empty paragraph right before the page-break paragraph)
# must stay untouched.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $text = $d.Paragraphs.Item($i).Range.Text
    if ($text -match "LOB1049: Estat") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    # Delete the three trailing paragraphs in reverse order so earlier
    # deletions don't shift the indices of the ones still to be removed.
    $d.Paragraphs.Item($targetIndex + 3).Range.Delete()
    $d.Paragraphs.Item($targetIndex + 2).Range.Delete()
    $d.Paragraphs.Item($targetIndex + 1).Range.Delete()
}
